# Generate Report for Handback
# Rows for 984b07f3-... now have handback data (mirroring the other two
# files which were already "Handed back: in sync with en-US"), and all
# three sheets are re-sorted by source file name:
#   984b07f3-...  <  ffff0f0c68c3-...  <  ffffff975c7fb3-...

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$ov.Range("B2").Value = "e2e\984b07f3-8c45-4712-a89f-06216d9f9206.md"
$ov.Range("G2").Value = "2016-08-16 11:06:44"

$ov.Range("A3").Value = "ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md"
$ov.Range("B3").Value = "e2e\ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md"

$ov.Range("A4").Value = "ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md"
$ov.Range("B4").Value = "e2e\ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md"
$ov.Range("E4").Value = "Handed back: in sync with en-US"
$ov.Range("F4").Value = "Handed back: in sync with en-US"
$ov.Range("G4").Value = "2016-08-16 11:03:35"

$ov.Hyperlinks.Item(1).Range.Value = "e2e\984b07f3-8c45-4712-a89f-06216d9f9206.md"
$ov.Hyperlinks.Item(1).TextToDisplay = "e2e\984b07f3-8c45-4712-a89f-06216d9f9206.md"
$ov.Hyperlinks.Item(2).Range.Value = "e2e\ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md"
$ov.Hyperlinks.Item(2).TextToDisplay = "e2e\ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md"
$ov.Hyperlinks.Item(3).Range.Value = "e2e\ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md"
$ov.Hyperlinks.Item(3).TextToDisplay = "e2e\ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$zh.Range("G2").Value = "984b07f3-8c45-4712-a89f-06216d9f9206.0a683d6ce457ecb89daf73c135c76f920d7d20cd.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-16 11:06:37"
$zh.Range("I2").Value = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$zh.Range("J2").Value = "984b07f3-8c45-4712-a89f-06216d9f9206.0a683d6ce457ecb89daf73c135c76f920d7d20cd.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-16 11:07:12"

$zh.Range("A3").Value = "ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md"
$zh.Range("F3").Value = "'False"

$zh.Range("A4").Value = "ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md"
$zh.Range("C4").Value = "Handed back: in sync with en-US"
$zh.Range("F4").Value = "'True"
$zh.Range("G4").Value = "885e7002-3dba-40ab-a7d3-33d242224785.0271c48b8eceb70fe07976a53047e849de7936c1.zh-cn.xlf"
$zh.Range("H4").Value = "2016-08-16 11:03:29"
$zh.Range("I4").Value = "885e7002-3dba-40ab-a7d3-33d242224785.md"
$zh.Range("J4").Value = "885e7002-3dba-40ab-a7d3-33d242224785.0271c48b8eceb70fe07976a53047e849de7936c1.zh-cn.xlf"
$zh.Range("K4").Value = "2016-08-16 11:03:57"
$zh.Range("P4").Value = ""

$zh.Hyperlinks.Item(1).Range.Value = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$zh.Hyperlinks.Item(1).TextToDisplay = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$zh.Hyperlinks.Item(2).Range.Value = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$zh.Hyperlinks.Item(2).TextToDisplay = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$zh.Hyperlinks.Item(3).Range.Value = "ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md"
$zh.Hyperlinks.Item(3).TextToDisplay = "ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md"
$zh.Hyperlinks.Item(4).Range.Value = "885e7002-3dba-40ab-a7d3-33d242224785.md"
$zh.Hyperlinks.Item(4).TextToDisplay = "885e7002-3dba-40ab-a7d3-33d242224785.md"
$zh.Hyperlinks.Item(5).Range.Value = "ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md"
$zh.Hyperlinks.Item(5).TextToDisplay = "ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md"
$zh.Hyperlinks.Item(6).Range.Value = "885e7002-3dba-40ab-a7d3-33d242224785.md"
$zh.Hyperlinks.Item(6).TextToDisplay = "885e7002-3dba-40ab-a7d3-33d242224785.md"

$zh.Columns.Item(16).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$de.Range("G2").Value = "984b07f3-8c45-4712-a89f-06216d9f9206.0a683d6ce457ecb89daf73c135c76f920d7d20cd.de-de.xlf"
$de.Range("H2").Value = "2016-08-16 11:06:44"
$de.Range("I2").Value = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$de.Range("J2").Value = "984b07f3-8c45-4712-a89f-06216d9f9206.0a683d6ce457ecb89daf73c135c76f920d7d20cd.de-de.xlf"
$de.Range("K2").Value = "2016-08-16 11:07:20"

$de.Range("A3").Value = "ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md"
$de.Range("F3").Value = "'False"

$de.Range("A4").Value = "ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md"
$de.Range("C4").Value = "Handed back: in sync with en-US"
$de.Range("F4").Value = "'True"
$de.Range("G4").Value = "885e7002-3dba-40ab-a7d3-33d242224785.0271c48b8eceb70fe07976a53047e849de7936c1.de-de.xlf"
$de.Range("H4").Value = "2016-08-16 11:03:35"
$de.Range("I4").Value = "885e7002-3dba-40ab-a7d3-33d242224785.md"
$de.Range("J4").Value = "885e7002-3dba-40ab-a7d3-33d242224785.0271c48b8eceb70fe07976a53047e849de7936c1.de-de.xlf"
$de.Range("K4").Value = "2016-08-16 11:04:12"
$de.Range("P4").Value = ""

$de.Hyperlinks.Item(1).Range.Value = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$de.Hyperlinks.Item(1).TextToDisplay = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$de.Hyperlinks.Item(2).Range.Value = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$de.Hyperlinks.Item(2).TextToDisplay = "984b07f3-8c45-4712-a89f-06216d9f9206.md"
$de.Hyperlinks.Item(3).Range.Value = "ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md"
$de.Hyperlinks.Item(3).TextToDisplay = "ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d.md"
$de.Hyperlinks.Item(4).Range.Value = "885e7002-3dba-40ab-a7d3-33d242224785.md"
$de.Hyperlinks.Item(4).TextToDisplay = "885e7002-3dba-40ab-a7d3-33d242224785.md"
$de.Hyperlinks.Item(5).Range.Value = "ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md"
$de.Hyperlinks.Item(5).TextToDisplay = "ffffff975c7fb3-1896-41b7-bfb4-207e2b9138d6.md"
$de.Hyperlinks.Item(6).Range.Value = "885e7002-3dba-40ab-a7d3-33d242224785.md"
$de.Hyperlinks.Item(6).TextToDisplay = "885e7002-3dba-40ab-a7d3-33d242224785.md"

$de.Columns.Item(16).AutoFit() | Out-Null
